$d = $word.ActiveDocument

# ------------------------------------------------------------------
# "Bajar datos de renta" -- fix the malformed number in ARTICULO SEXTO:
#   "Fijase en $a. 3. 000 ..."  ->  "Fijase en $a 3.000 ..."
# i.e. drop the stray period after "$a" (keep the space that follows
# it) and drop the stray space between "3." and "000" so the figure
# reads as a normal "3.000" amount. The editor's cursor ends up right
# after the new "3." -- which is where Word's automatic "_GoBack"
# bookmark is left, so we move that bookmark there too.
# ------------------------------------------------------------------

# 1. Drop the old "_GoBack" bookmark whereever Word last left it
#    (it currently sits at the very end of the document, after the
#    closing "ARCHIVESE.").
$goBack = $null
foreach ($dummy in 1..1) {
    $goBack = $d.Bookmarks.Item("_GoBack")
}
if ($goBack -ne $null) {
    $goBack.Delete()
}

# 2. "3. 000" -> "3.000": remove the space between the period and the
#    thousands group.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found1 = $rng.Find.Execute('3. 000', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found1) {
    # Narrow to just the space that sits right after "3." (index 2
    # inside the found match: "3", ".", " ").
    $periodEnd = $rng.Start + 2
    $spaceRng = $d.Range($periodEnd, $periodEnd + 1)
    $spaceRng.Text = ""

    # 3. Re-create "_GoBack" right after the new "3." -- between the
    #    period and "000" -- matching where Word leaves the cursor
    #    after the edit.
    $d.Bookmarks.Add("_GoBack", $d.Range($periodEnd, $periodEnd))
}

# 4. "$a. 3" -> "$a 3": remove the period right after "$a", leaving
#    the space that was already after it.
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute('$a. 3', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found2) {
    $periodPos = $rng2.Start + 2
    $dotRng = $d.Range($periodPos, $periodPos + 1)
    $dotRng.Text = ""
}
